# Update the theme/party counts with the new lemmatized values, merge the
# "Tourisme" row into the updated "Économie"/"Éducation" rows, and drop the
# now-redundant last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-12 (themes unchanged, only numeric columns change).
$data = @(
    @{ Row = 2;  B = 5;  C = 14;  D = 1;  E = 7 }   # Agriculture
    @{ Row = 3;  B = 2;  C = 10;  D = 0;  E = 12 }  # Droits_Femme
    @{ Row = 4;  B = 23; C = 19;  D = 6;  E = 24 }  # Emploi
    @{ Row = 5;  B = 10; C = 79;  D = 10; E = 3 }   # Environnement
    @{ Row = 6;  B = 22; C = 56;  D = 24; E = 23 }  # Gouvernance
    @{ Row = 7;  B = 0;  C = 19;  D = 0;  E = 3 }   # Infrastructure
    @{ Row = 8;  B = 4;  C = 0;   D = 6;  E = 0 }   # Jeunesse
    @{ Row = 9;  B = 14; C = 62;  D = 9;  E = 18 }  # Justice
    @{ Row = 10; B = 6;  C = 7;   D = 0;  E = 2 }   # Logement
    @{ Row = 11; B = 12; C = 10;  D = 0;  E = 17 }  # Santé
    @{ Row = 12; B = 15; C = 44;  D = 5;  E = 38 }  # Social
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 2).Value = $entry.B
    $ws.Cells.Item($r, 3).Value = $entry.C
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 5).Value = $entry.E
}

# Row 13 used to be "Tourisme"; it is replaced with updated "Économie" data.
$ws.Cells.Item(13, 1).Value = "Économie"
$ws.Cells.Item(13, 2).Value = 13
$ws.Cells.Item(13, 3).Value = 116
$ws.Cells.Item(13, 4).Value = 2
$ws.Cells.Item(13, 5).Value = 27

# Row 14 used to be "Économie"; it is replaced with updated "Éducation" data.
$ws.Cells.Item(14, 1).Value = "Éducation"
$ws.Cells.Item(14, 2).Value = 10
$ws.Cells.Item(14, 3).Value = 9
$ws.Cells.Item(14, 4).Value = 4
$ws.Cells.Item(14, 5).Value = 16

# Row 15 (old "Éducation") is now merged into row 14 above, so delete it.
$ws.Rows.Item(15).Delete()
